$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values stay text (matching source data which uses
# strings like "25.691.89" / "0.4931" rather than numbers), regardless of
# whether Excel would otherwise auto-detect them as numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.691.89'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").Value = '1.742.84'
$ws.Range("E3").Value = '  -5.39%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '238.43'
$ws.Range("E5").Value = '  -8.40%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = '0.4931'
$ws.Range("E7").Value = '  -6.46%  '
$ws.Range("D8").Value = '41.42'
$ws.Range("E8").Value = '  -7.78%  '
$ws.Range("D9").Value = '0.2448'
$ws.Range("E9").Value = '  -22.60%  '
$ws.Range("D10").Value = '0.05954'
$ws.Range("E10").Value = '  -12.35%  '
$ws.Range("D11").Value = '1.745.08'
$ws.Range("E11").Value = '  -5.08%  '
$ws.Range("D12").Value = '0.06778'
$ws.Range("E12").Value = '  -12.97%  '
$ws.Range("D13").Value = '14.76'
$ws.Range("E13").Value = '  -22.73%  '
$ws.Range("D14").Value = '4.467'
$ws.Range("E14").Value = '  -10.89%  '
$ws.Range("D15").Value = '77.01'
$ws.Range("E15").Value = '  -12.74%  '
$ws.Range("D16").Value = '0.5804'
$ws.Range("E16").Value = '  -25.85%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '25.740.71'
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("D20").Value = '11.50'
$ws.Range("E20").Value = '  -17.27%  '
$ws.Range("D21").Value = '0.000006456'
$ws.Range("E21").Value = '  -18.50%  '
$ws.Range("D22").Value = '1.965.76'
$ws.Range("E22").Value = '  -5.20%  '
$ws.Range("D23").Value = '3.973'
$ws.Range("E23").Value = '  -13.83%  '
$ws.Range("D24").Value = '7.883'
$ws.Range("E24").Value = '  -15.56%  '
$ws.Range("D25").Value = '5.010'
$ws.Range("E25").Value = '  -16.40%  '
$ws.Range("D26").Value = '135.86'
$ws.Range("E26").Value = '  -4.97%  '
$ws.Range("D27").Value = '1.494'
$ws.Range("E27").Value = '  -11.15%  '
$ws.Range("D28").Value = '1.830'
$ws.Range("E28").Value = '  -17.68%  '
$ws.Range("D29").Value = '14.54'
$ws.Range("E29").Value = '  -14.56%  '
$ws.Range("D30").Value = '100.71'
$ws.Range("E30").Value = '  -9.16%  '
$ws.Range("D31").Value = '3.794'
$ws.Range("E31").Value = '  -9.86%  '
$ws.Range("D32").Value = '0.08096'
$ws.Range("E32").Value = '  -6.84%  '
$ws.Range("D33").Value = '3.337'
$ws.Range("E33").Value = '  -18.15%  '
$ws.Range("D34").Value = '0.04407'
$ws.Range("E34").Value = '  -9.28%  '
$ws.Range("D35").Value = '0.9997'
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").Value = '2.650'
$ws.Range("E36").Value = '  -7.53%  '
$ws.Range("D37").Value = '1.015'
$ws.Range("E37").Value = '  -10.83%  '
$ws.Range("D38").Value = '0.6055'
$ws.Range("E38").Value = '  -17.32%  '
$ws.Range("D39").Value = '2.704'
$ws.Range("E39").Value = '  -12.68%  '
$ws.Range("D40").Value = '2.051'
$ws.Range("E40").Value = '  -12.54%  '
$ws.Range("D41").Value = '1.000'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '103.54'
$ws.Range("E42").Value = '  -5.21%  '
$ws.Range("D43").Value = '0.01500'
$ws.Range("E43").Value = '  -13.41%  '
$ws.Range("D44").Value = '0.7776'
$ws.Range("E44").Value = '  -13.99%  '
$ws.Range("D45").Value = '5.181'
$ws.Range("E45").Value = '  -12.31%  '
$ws.Range("D46").Value = '0.3758'
$ws.Range("E46").Value = '  -21.96%  '
$ws.Range("D47").Value = '0.05117'
$ws.Range("E47").Value = '  -12.13%  '
$ws.Range("D48").Value = '0.1076'
$ws.Range("E48").Value = '  -13.47%  '
$ws.Range("D49").Value = '5.946'
$ws.Range("E49").Value = '  -22.79%  '
$ws.Range("D50").Value = '30.22'
$ws.Range("E50").Value = '  -13.38%  '
$ws.Range("D51").Value = '52.58'
$ws.Range("E51").Value = '  -12.44%  '
